$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26: Student milestone changed from "I" to "II"
$ws.Range("E26").Value = "II"

# Rows 29-31: mark Milestone II complete (student picked Milestone II, with "X" completed)
$ws.Range("E29").Value = "II"
$ws.Range("F29").Value = "X"

$ws.Range("E30").Value = "II"
$ws.Range("F30").Value = "X"

$ws.Range("E31").Value = "II"
$ws.Range("F31").Value = "X"

# Rows 46-47: mark Milestone II complete
$ws.Range("E46").Value = "II"
$ws.Range("F46").Value = "X"

$ws.Range("E47").Value = "II"
$ws.Range("F47").Value = "X"

# Rows 73-74 (Milestone II Complete(X) summary column D)
$ws.Range("D73").Value = "X"
$ws.Range("D74").Value = "X"

# Row 79: add new source citation link
$ws.Range("A79").Value = "https://www.models-resource.com/pc_computer/zootycoon2/model/17846/"

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E10").Select()
